$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.224.93"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "2.386.92"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'562.96"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("D6").Value = "'138.24"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.537"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "2.387.40"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").Value = "'0.105"
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("D11").Value = "'0.160"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "'5.05"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").Value = "'25.64"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("D17").Value = "60.135.90"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'8.42"
$ws.Range("E18").Value = "  +16.03%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "2.389.49"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").Value = "'10.60"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "'324.82"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "'4.03"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").Value = "'6.08"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  -8.25%  "
$ws.Range("D26").Value = "'64.40"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").Value = "'554.51"
$ws.Range("E27").Value = "  -4.83%  "
$ws.Range("D28").Value = "'7.94"
$ws.Range("E28").Value = "  -13.17%  "
$ws.Range("D29").Value = "2.522.51"
$ws.Range("D30").Value = "0.0₃0896"
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").Value = "'7.91"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "'1.29"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").Value = "'1.80"
$ws.Range("E33").Value = "  -3.74%  "
$ws.Range("D34").Value = "'0.130"
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").Value = "'153.51"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("D37").Value = "'1.41"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'0.366"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").Value = "'4.50"
$ws.Range("E39").Value = "  -4.52%  "
$ws.Range("D40").Value = "'18.24"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").Value = "'5.02"
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'41.08"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.63"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.27"
$ws.Range("E45").Value = "  -4.34%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0279"
$ws.Range("E46").Value = "  -5.33%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'142.95"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "'3.49"
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.584"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").Value = "'0.0498"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'18.81"
$ws.Range("E51").Value = "  -4.33%  "
